$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy a cell with the same style (s="11") to use as a pure-format source for PasteSpecial
$ws.Range("C9").Copy()

$ws.Range("D9").Value = "'0.2343"
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("F9").Value = "'0.5932"
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("D11").Value = "'0.4896"
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("F11").Value = "'0.3410"
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("D13").Value = "'0.2376"
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("F13").Value = "'0.6916"
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0.5937"
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("F15").Value = "'0.3518"
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("D17").Value = "'0.6982"
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("F17").Value = "'0.0978"
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("D19").Value = "'0.4604"
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("F19").Value = "'0.3871"
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("D21").Value = "'0.0272"
$ws.Range("D21").PasteSpecial(-4122)
$ws.Range("F21").Value = "'0.9546"
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("D23").Value = "'0.1985"
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("F23").Value = "'0.3916"
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("D25").Value = "'0.8158"
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("F25").Value = "'0.5864"
$ws.Range("F25").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0.0804"
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("F27").Value = "'0.9526"
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("D29").Value = "'0.4208"
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("F29").Value = "'0.0379"
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("D31").Value = "'0.6972"
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("F31").Value = "'0.1794"
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("D33").Value = "'0.2505"
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("F33").Value = "'0.1763"
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("D35").Value = "'0.6450"
$ws.Range("D35").PasteSpecial(-4122)
$ws.Range("F35").Value = "'0.5890"
$ws.Range("F35").PasteSpecial(-4122)
$ws.Range("D37").Value = "'0.7140"
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("F37").Value = "'0.7115"
$ws.Range("F37").PasteSpecial(-4122)
$ws.Range("D39").Value = "'0.2594"
$ws.Range("D39").PasteSpecial(-4122)
$ws.Range("F39").Value = "'0.9770"
$ws.Range("F39").PasteSpecial(-4122)
$ws.Range("D41").Value = "'0.6945"
$ws.Range("D41").PasteSpecial(-4122)
$ws.Range("F41").Value = "'0.7974"
$ws.Range("F41").PasteSpecial(-4122)
$ws.Range("D43").Value = "'0.7965"
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("F43").Value = "'0.9127"
$ws.Range("F43").PasteSpecial(-4122)
$ws.Range("D45").Value = "'0.6083"
$ws.Range("D45").PasteSpecial(-4122)
$ws.Range("F45").Value = "'0.5479"
$ws.Range("F45").PasteSpecial(-4122)
$ws.Range("D47").Value = "'0.0027"
$ws.Range("D47").PasteSpecial(-4122)
$ws.Range("F47").Value = "'0.2314"
$ws.Range("F47").PasteSpecial(-4122)
$ws.Range("D49").Value = "'0.3293"
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("F49").Value = "'0.5688"
$ws.Range("F49").PasteSpecial(-4122)
$ws.Range("D51").Value = "'0.3072"
$ws.Range("D51").PasteSpecial(-4122)
$ws.Range("F51").Value = "'0.2269"
$ws.Range("F51").PasteSpecial(-4122)
$ws.Range("D53").Value = "'0.8420"
$ws.Range("D53").PasteSpecial(-4122)
$ws.Range("F53").Value = "'0.0385"
$ws.Range("F53").PasteSpecial(-4122)
$ws.Range("D55").Value = "'0.0715"
$ws.Range("D55").PasteSpecial(-4122)
$ws.Range("F55").Value = "'0.4508"
$ws.Range("F55").PasteSpecial(-4122)
$ws.Range("D57").Value = "'0.1284"
$ws.Range("D57").PasteSpecial(-4122)
$ws.Range("F57").Value = "'0.0176"
$ws.Range("F57").PasteSpecial(-4122)
$ws.Range("D59").Value = "'0.6429"
$ws.Range("D59").PasteSpecial(-4122)
$ws.Range("F59").Value = "'0.5555"
$ws.Range("F59").PasteSpecial(-4122)
$ws.Range("D61").Value = "'0.7460"
$ws.Range("D61").PasteSpecial(-4122)
$ws.Range("F61").Value = "'0.4122"
$ws.Range("F61").PasteSpecial(-4122)
$ws.Range("D63").Value = "'0.6371"
$ws.Range("D63").PasteSpecial(-4122)
$ws.Range("F63").Value = "'0.2979"
$ws.Range("F63").PasteSpecial(-4122)
